# Updated to simulate to 9-feb
# Insert a new column before column H ("Cotton.Leaf.SpecificArea"),
# shifting all existing H..O columns right by one (H..P), and populate
# the new column with =F{row}/M{row} formulas for the rows where new
# leaf-area data became available. Also add the new underlying Wt data
# point for row 18 (Cotton.Leaf.Wt / Cotton.Stem.Wt) that the formula
# depends on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert the new column. Excel automatically shifts cell contents,
#    formulas, styles and column width metadata one column to the right.
$ws.Columns("H").Insert()

# 2. Header for the newly inserted column.
$ws.Range("H1").Value = "Cotton.Leaf.SpecificArea"

# 3. New formula cells: SpecificArea = LAI / LeafWt (F / M after the shift).
$ws.Range("H12").Formula = "=F12/M12"
$ws.Range("H15").Formula = "=F15/M15"
$ws.Range("H18").Formula = "=F18/M18"
$ws.Range("H27").Formula = "=F27/M27"
$ws.Range("H29").Formula = "=F29/M29"
$ws.Range("H32").Formula = "=F32/M32"

# 4. New underlying data point added for row 18 (Cotton.Leaf.Wt / Cotton.Stem.Wt),
#    needed for the new SpecificArea formula in H18. N18 keeps the same
#    "0.00" number format already used by its column neighbours (K18/L18).
$ws.Range("M18").Value = 282.98
$ws.Range("N18").Value = 425.3
$ws.Range("N18").NumberFormat = "0.00"

# 5. Keep the worksheet's filter database defined name in sync with the
#    widened data range.
$wb.Names.Item("CottonObserved!_FilterDatabase").RefersTo = "=CottonObserved!`$A`$1:`$EQ`$2578"
